$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Copy number formats/styles from column F (the old column D, now shifted) into new D:E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D and E values for each data row
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 395300
$ws.Range("E8").Value2 = 406800
$ws.Range("D9").Value2 = 297000
$ws.Range("E9").Value2 = 308000
$ws.Range("D10").Value2 = 98300
$ws.Range("E10").Value2 = 98800
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = "NA"
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 373200
$ws.Range("E17").Value2 = 384400
$ws.Range("D18").Value2 = 22100
$ws.Range("E18").Value2 = 22400
$ws.Range("D20").Value2 = -13100
$ws.Range("E20").Value2 = -12700
$ws.Range("D21").Value2 = 16000
$ws.Range("E21").Value2 = 17000
$ws.Range("D22").Value2 = 0
$ws.Range("E22").Value2 = 0
$ws.Range("D23").Value2 = 8900
$ws.Range("E23").Value2 = 9600
$ws.Range("D24").Value2 = 37700
$ws.Range("E24").Value2 = 400
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = -28700
$ws.Range("E26").Value2 = 9300
$ws.Range("D27").Value2 = -28700
$ws.Range("E27").Value2 = 9300
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 35000
$ws.Range("E29").Value2 = -2000
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 13100
$ws.Range("E32").Value2 = 12700
$ws.Range("D33").Value2 = 6300
$ws.Range("E33").Value2 = 7300
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 6300
$ws.Range("E35").Value2 = 7300
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 25200
$ws.Range("E41").Value2 = 46200
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 299600
$ws.Range("E43").Value2 = 285800
$ws.Range("D44").Value2 = 912700
$ws.Range("E44").Value2 = 884200
$ws.Range("D45").Value2 = 19300
$ws.Range("E45").Value2 = 15300
$ws.Range("D46").Value2 = 1256700
$ws.Range("E46").Value2 = 1231500
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 44000
$ws.Range("E48").Value2 = 44200
$ws.Range("D49").Value2 = 426300
$ws.Range("E49").Value2 = 430100
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 86000
$ws.Range("E52").Value2 = 83700
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 1813100
$ws.Range("E54").Value2 = 1789500
$ws.Range("D57").Value2 = 184900
$ws.Range("E57").Value2 = 180500
$ws.Range("D58").Value2 = 96200
$ws.Range("E58").Value2 = 76200
$ws.Range("D59").Value2 = 41200
$ws.Range("E59").Value2 = 45100
$ws.Range("D60").Value2 = 322300
$ws.Range("E60").Value2 = 301800
$ws.Range("D61").Value2 = 769800
$ws.Range("E61").Value2 = 774100
$ws.Range("D62").Value2 = 22700
$ws.Range("E62").Value2 = 21100
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 1114800
$ws.Range("E66").Value2 = 1097000
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 337100
$ws.Range("E72").Value2 = 330800
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 698300
$ws.Range("E76").Value2 = 692500
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 6300
$ws.Range("E81").Value2 = 7300
$ws.Range("D83").Value2 = 7100
$ws.Range("E83").Value2 = 7300
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = -32400
$ws.Range("E89").Value2 = 36900
$ws.Range("D91").Value2 = -2200
$ws.Range("E91").Value2 = -1700
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -2200
$ws.Range("E94").Value2 = -1700
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = 13800
$ws.Range("E100").Value2 = -34500
$ws.Range("D101").Value2 = -200
$ws.Range("E101").Value2 = -200
$ws.Range("D102").Value2 = -21000
$ws.Range("E102").Value2 = 600

# Restated historical values for a few rows (F,G,H corrected in addition to new D,E)
$ws.Range("F24").Value2 = 5900
$ws.Range("G24").Value2 = 6200
$ws.Range("H24").Value2 = 4300
$ws.Range("F26").Value2 = 11000
$ws.Range("G26").Value2 = 14900
$ws.Range("H26").Value2 = 8700
$ws.Range("F27").Value2 = 11000
$ws.Range("G27").Value2 = 14900
$ws.Range("H27").Value2 = 8700
$ws.Range("F29").Value2 = -200
$ws.Range("G29").Value2 = 100
$ws.Range("H29").Value2 = -9100